# Insert a new row at row 218 on the active sheet, pushing the existing
# rows 218:257 down to 219:258 (dimension grows from A1:R257 to A1:R258),
# then populate the newly-inserted row 218 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 218:257 down by one to make room for the new record.
$ws.Rows.Item(218).Insert()

# Fill in the new row 218 with the new data point.
$ws.Cells.Item(218, 1).Value = 5
$ws.Cells.Item(218, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(218, 3).Value = "Maule"
$ws.Cells.Item(218, 4).Value = 44504
$ws.Cells.Item(218, 5).Value = 7
$ws.Cells.Item(218, 6).Value = 100112043
$ws.Cells.Item(218, 7).Value = "Pepino ensalada"
$ws.Cells.Item(218, 8).Value = "Sin especificar"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 400
$ws.Cells.Item(218, 11).Value = 6000
$ws.Cells.Item(218, 12).Value = 6000
$ws.Cells.Item(218, 13).Value = 6000
$ws.Cells.Item(218, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(218, 15).Value = "Región del Maule"
$ws.Cells.Item(218, 16).Value = 75
$ws.Cells.Item(218, 17).Value = 80
$ws.Cells.Item(218, 18).Value = "Hortaliza"
